$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.748.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.274.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.97%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.40%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'230.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.85%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'63.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.15%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.77%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +8.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'57.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'25.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +14.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.614.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.08%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.15%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.816"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.296.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.77%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.653.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.36%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0951"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'73.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.83%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'250.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +4.81%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.47%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'171.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.24%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.52%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'20.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.21%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0692"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.39%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.21%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.12%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0247"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.76%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'TerraClassic"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.000222"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -10.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Celestia"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'10.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +23.70%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.95%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.82%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.08%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.31%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'97.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.58%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.482.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.63%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'16.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.73%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.13%  "
$ws.Range("E51").Style = "Normal"
